# Automatische test-sync: 2025-06-19 16:00:10
# Appends the latest incoming mail to the "Logs" sheet and refreshes the
# "Dashboard" category summary accordingly.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append the new mail entry as row 17 ---------------------
$newRow = 17

$logs.Range("A" + $newRow).Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B" + $newRow).Value = "mailmind.test@zohomail.eu"
$logs.Range("C" + $newRow).Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D" + $newRow).Value = "Bestelling"
$logs.Range("F" + $newRow).Value = "2025-06-19 15:58:10"
$logs.Range("G" + $newRow).Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ------
$catFormats = $logs.Range("D2:D16").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D" + $newRow))
}

$answeredFormats = $logs.Range("G2:G16").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G" + $newRow))
}

# --- Dashboard sheet: refresh the category counts --------------------------
# "Bestelling" now has 2 occurrences and moves above "Informatieaanvraag" (1)
$dashboard.Range("A5").Value = "Bestelling"
$dashboard.Range("B5").Value = 2
$dashboard.Range("A6").Value = "Informatieaanvraag"
$dashboard.Range("B6").Value = 1
